$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 511, shifting rows 511:532
# down to 512:533 (and expanding the used range to A1:R533).
$ws.Rows.Item(511).Insert()

# Populate the newly-inserted row 511 with the new weekly price-report data.
$ws.Range("A511").Value = 6
$ws.Range("B511").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C511").Value = "Metropolitana"
$ws.Range("D511").Value = 44753
$ws.Range("E511").Value = 13
$ws.Range("F511").Value = 100112044
$ws.Range("G511").Value = "Perejil"
$ws.Range("H511").Value = "Sin especificar"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 80
$ws.Range("K511").Value = 19000
$ws.Range("L511").Value = 20000
$ws.Range("M511").Value = 19375
$ws.Range("N511").Value = "$/docena de atados"
$ws.Range("O511").Value = "Región Metropolitana"
$ws.Range("P511").Value = 6458
$ws.Range("Q511").Value = 3
$ws.Range("R511").Value = "Hortaliza"
